# Update column G ("K") values per regenerated save_data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 2
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 2
    19 = 1
    20 = 1
    21 = 0
    22 = 0
    23 = 2
    24 = 1
    25 = 2
    26 = 2
    28 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
